# Updated cryptos list on Thu Feb 23 06:11:39 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "24.418.50"
$ws.Range("E2").Value = "  +1.18%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "1.667.10"
$ws.Range("E3").Value = "  +1.39%  "

# Row 4 (TetherUSD)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.99"
$ws.Range("E5").Value = "  +1.50%  "

# Row 6 (USDC)
$ws.Range("E6").Value = "  -0.03%  "

# Row 7 (XRP)
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3960"
$ws.Range("E7").Value = "  +1.08%  "

# Row 8 (Cardano)
$ws.Range("E8").Value = "  +1.19%  "

# Row 9 (OKB)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.42"
$ws.Range("E9").Value = "  +6.83%  "

# Row 10 (Polygon)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.401"
$ws.Range("E10").Value = "  +2.94%  "

# Row 11 (BinanceUSD)
$ws.Range("E11").Value = "  -0.02%  "

# Row 12 (Dogecoin)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08591"
$ws.Range("E12").Value = "  +1.33%  "

# Row 13 (Solana)
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.42"
$ws.Range("E13").Value = "  +1.28%  "

# Row 14 (Polkadot)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.335"
$ws.Range("E14").Value = "  +2.48%  "

# Row 15 (ShibaInu)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001353"
$ws.Range("E15").Value = "  +5.01%  "

# Row 16 (Chainlink)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.897"
$ws.Range("E16").Value = "  +5.24%  "

# Row 17 (WrappedEther)
$ws.Range("D17").Value = "1.664.35"
$ws.Range("E17").Value = "  +1.06%  "

# Row 18 (Litecoin)
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.31"
$ws.Range("E18").Value = "  +1.12%  "

# Row 19 (TRON)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06996"
$ws.Range("E19").Value = "  +0.60%  "

# Row 20 (Avalanche)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.59"
$ws.Range("E20").Value = "  -1.69%  "

# Row 21 (Uniswap)
$ws.Range("E21").Value = "  +0.60%  "

# Row 22 (Dai)
$ws.Range("E22").Value = "  -0.08%  "

# Row 23 (Cosmos)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.77"
$ws.Range("E23").Value = "  +0.18%  "

# Row 24 (WrappedBTC)
$ws.Range("D24").Value = "24.421.49"
$ws.Range("E24").Value = "  +1.12%  "

# Row 25 (Toncoin -> LidoDAOToken)
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.075"
$ws.Range("E25").Value = "  +12.42%  "

# Row 26 (LidoDAOToken -> Toncoin)
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.420"
$ws.Range("E26").Value = "  +2.78%  "

# Row 27 (EthereumClassic)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.57"
$ws.Range("E27").Value = "  +0.18%  "

# Row 28 (Monero)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.35"
$ws.Range("E28").Value = "  -0.44%  "

# Row 29 (BitcoinCash)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "142.83"
$ws.Range("E29").Value = "  +0.81%  "

# Row 30 (HuobiToken)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.413"
$ws.Range("E30").Value = "  +1.03%  "

# Row 31 (Filecoin)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.096"
$ws.Range("E31").Value = "  -8.80%  "

# Row 32 (WEMIXTOKEN)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.532"
$ws.Range("E32").Value = "  +1.97%  "

# Row 33 (WrappedliquidstakedEther2.0)
$ws.Range("D33").Value = "1.846.45"
$ws.Range("E33").Value = "  +1.14%  "

# Row 34 (ImmutableX)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.069"
$ws.Range("E34").Value = "  +8.75%  "

# Row 35 (Hedera)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08266"
$ws.Range("E35").Value = "  +2.09%  "

# Row 36 (VeChain)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03045"
$ws.Range("E36").Value = "  +3.62%  "

# Row 37 (InternetComputer(DFINITY))
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.902"
$ws.Range("E37").Value = "  -3.94%  "

# Row 38 (Algorand -> FraxShare)
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.14"
$ws.Range("E38").Value = "  +10.33%  "

# Row 39 (FraxShare -> Algorand)
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2767"
$ws.Range("E39").Value = "  +2.14%  "

# Row 40 (Stellar)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09227"
$ws.Range("E40").Value = "  -0.27%  "

# Row 41 (TheSandbox)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7747"
$ws.Range("E41").Value = "  +1.13%  "

# Row 42 (Aptos)
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.82"
$ws.Range("E42").Value = "  +4.77%  "

# Row 43 (TrustWalletToken)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.447"
$ws.Range("E43").Value = "  -2.11%  "

# Row 44 (EnergySwap)
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.66"
$ws.Range("E44").Value = "  +2.91%  "

# Row 45 (Decentraland)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7129"
$ws.Range("E45").Value = "  +3.22%  "

# Row 46 (NEARProtocol)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.540"
$ws.Range("E46").Value = "  +1.98%  "

# Row 47 (PancakeSwap)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.133"
$ws.Range("E47").Value = "  +1.06%  "

# Row 48 (Frax)
$ws.Range("E48").Value = "  -0.04%  "

# Row 49 (Cronos)
$ws.Range("E49").Value = "  +0.37%  "

# Row 50 (Quant)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.63"
$ws.Range("E50").Value = "  +1.66%  "

# Row 51 (Flow)
$ws.Range("E51").Value = "  +0.22%  "
